# Auto-generated edit script for "Vega Modelo de Temuco - Puerro" sheet.
# A new daily price record is inserted at the top of this product block
# (the rows immediately following become one row lower), so row 179 receives
# brand-new data and every following row (180-270) takes on the values that
# used to live one row above it. The row that used to be last (270) is
# copied down into a brand-new row 271, extending the sheet dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> column index, for the fields that actually change (D,I,J,K,L,M,O,P)
$colIndex = @{ "D" = 4; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13; "O" = 15; "P" = 16 }

# Each entry: row number, then the new values for columns D, I, J, K, L, M, O, P
$updates = @(
    @(179, 45001, 'Primera', 80, 14000, 14000, 14000, 'Provincia de Cautín', 1167),
    @(180, 44463, 'Primera', 30, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(181, 44391, 'Primera', 60, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(182, 44812, 'Primera', 40, 17000, 17000, 17000, 'Provincia de Cautín', 1417),
    @(183, 44847, 'Primera', 95, 18000, 18000, 18000, 'Provincia de Cautín', 1500),
    @(184, 44676, 'Primera', 30, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(185, 44889, 'Primera', 65, 18000, 18000, 18000, 'Provincia de Cautín', 1500),
    @(186, 44272, 'Primera', 20, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(187, 44309, 'Primera', 55, 6000, 7000, 6455, 'Provincia de Cautín', 538),
    @(188, 44314, 'Primera', 40, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(189, 44263, 'Primera', 80, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(190, 44410, 'Primera', 80, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(191, 44477, 'Primera', 50, 6000, 7000, 6600, 'Provincia de Cautín', 550),
    @(192, 44790, 'Primera', 65, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(193, 44330, 'Primera', 40, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(194, 44441, 'Primera', 50, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(195, 44756, 'Primera', 30, 17000, 17000, 17000, 'Provincia de Cautín', 1417),
    @(196, 44819, 'Primera', 30, 17000, 17000, 17000, 'Provincia de Cautín', 1417),
    @(197, 44841, 'Primera', 30, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(198, 44431, 'Primera', 65, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(199, 44942, 'Primera', 55, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(200, 44294, 'Primera', 65, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(201, 44855, 'Primera', 20, 1600, 1600, 1600, 'Provincia de Cautín', 133),
    @(202, 44519, 'Primera', 155, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(203, 44798, 'Primera', 30, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(204, 44722, 'Primera', 40, 15000, 16000, 15500, 'Provincia de Cautín', 1292),
    @(205, 44691, 'Primera', 55, 12000, 12000, 12000, 'Provincia de Cautín', 1000),
    @(206, 44420, 'Primera', 155, 7000, 8000, 7484, 'Provincia de Cautín', 624),
    @(207, 44637, 'Primera', 40, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(208, 44419, 'Primera', 55, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(209, 44186, 'Primera', 40, 15000, 15000, 15000, 'Provincia de Cautín', 1250),
    @(210, 44680, 'Primera', 40, 12000, 12000, 12000, 'Provincia de Cautín', 1000),
    @(211, 44222, 'Primera', 55, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(212, 44343, 'Primera', 50, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(213, 44433, 'Primera', 54, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(214, 44350, 'Primera', 65, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(215, 44721, 'Primera', 30, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(216, 44719, 'Primera', 40, 12000, 13000, 12500, 'Provincia de Cautín', 1042),
    @(217, 44399, 'Primera', 70, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(218, 44273, 'Primera', 90, 9000, 10000, 9444, 'Provincia de Cautín', 787),
    @(219, 44455, 'Primera', 40, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(220, 44764, 'Primera', 65, 17000, 17000, 17000, 'Provincia de Cautín', 1417),
    @(221, 44782, 'Primera', 20, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(222, 44515, 'Primera', 110, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(223, 44512, 'Primera', 40, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(224, 44504, 'Primera', 110, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(225, 44435, 'Primera', 299, 7000, 8000, 7900, 'Provincia de Cautín', 658),
    @(226, 44425, 'Primera', 30, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(227, 44987, 'Primera', 30, 14000, 14000, 14000, 'Provincia de Cautín', 1167),
    @(228, 44371, 'Primera', 30, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(229, 44392, 'Primera', 65, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(230, 44286, 'Primera', 80, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(231, 44827, 'Primera', 20, 15000, 15000, 15000, 'Provincia de Cautín', 1250),
    @(232, 44529, 'Primera', 65, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(233, 44411, 'Primera', 40, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(234, 44789, 'Primera', 65, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(235, 44383, 'Primera', 70, 8000, 9000, 8714, 'Provincia de Cautín', 726),
    @(236, 44238, 'Primera', 75, 13000, 14000, 13467, 'Provincia de Cautín', 1122),
    @(237, 44701, 'Primera', 50, 12000, 12000, 12000, 'Provincia de Cautín', 1000),
    @(238, 44508, 'Primera', 80, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(239, 44204, 'Primera', 30, 12000, 14000, 12667, 'Provincia de Cautín', 1056),
    @(240, 44663, 'Primera', 20, 12000, 12000, 12000, 'Provincia de Cautín', 1000),
    @(241, 44938, 'Primera', 65, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(242, 44760, 'Primera', 65, 15000, 15000, 15000, 'Provincia de Cautín', 1250),
    @(243, 44285, 'Primera', 20, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(244, 44291, 'Primera', 75, 8000, 10000, 8933, 'Provincia de Cautín', 744),
    @(245, 44292, 'Primera', 45, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(246, 44386, 'Primera', 70, 8000, 9000, 8429, 'Provincia de Cautín', 702),
    @(247, 44428, 'Primera', 30, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(248, 44868, 'Primera', 80, 18000, 19000, 18500, 'Provincia de Cautín', 1542),
    @(249, 44278, 'Primera', 45, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(250, 44795, 'Primera', 90, 17000, 18000, 17444, 'Provincia de Cautín', 1454),
    @(251, 44202, 'Primera', 20, 15000, 15000, 15000, 'Provincia de Cautín', 1250),
    @(252, 44610, 'Primera', 30, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(253, 44568, 'Primera', 20, 15000, 15000, 15000, 'Provincia de Cautín', 1250),
    @(254, 44473, 'Primera', 50, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(255, 44306, 'Primera', 35, 9000, 9000, 9000, 'Provincia de Cautín', 750),
    @(256, 44620, 'Primera', 10, 13000, 13000, 13000, 'Provincia de Cautín', 1083),
    @(257, 44413, 'Primera', 50, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(258, 44257, 'Primera', 50, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(259, 44426, 'Primera', 30, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(260, 44342, 'Primera', 30, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(261, 44606, 'Primera', 20, 10000, 10000, 10000, 'Provincia de Cautín', 833),
    @(262, 44567, 'Primera', 30, 15000, 15000, 15000, 'Región de La Araucanía', 1250),
    @(263, 44567, 'Segunda', 10, 13000, 13000, 13000, 'Región de La Araucanía', 1083),
    @(264, 44364, 'Primera', 65, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(265, 44215, 'Primera', 40, 12000, 13000, 12500, 'Provincia de Cautín', 1042),
    @(266, 44771, 'Primera', 20, 16000, 16000, 16000, 'Provincia de Cautín', 1333),
    @(267, 44498, 'Primera', 60, 7000, 7000, 7000, 'Provincia de Cautín', 583),
    @(268, 44988, 'Primera', 30, 14000, 14000, 14000, 'Provincia de Cautín', 1167),
    @(269, 44414, 'Primera', 40, 8000, 8000, 8000, 'Provincia de Cautín', 667),
    @(270, 44236, 'Primera', 55, 12000, 12000, 12000, 'Provincia de Cautín', 1000)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, $colIndex["D"]).Value = $u[1]
    $ws.Cells.Item($r, $colIndex["I"]).Value = $u[2]
    $ws.Cells.Item($r, $colIndex["J"]).Value = $u[3]
    $ws.Cells.Item($r, $colIndex["K"]).Value = $u[4]
    $ws.Cells.Item($r, $colIndex["L"]).Value = $u[5]
    $ws.Cells.Item($r, $colIndex["M"]).Value = $u[6]
    $ws.Cells.Item($r, $colIndex["O"]).Value = $u[7]
    $ws.Cells.Item($r, $colIndex["P"]).Value = $u[8]
}

# New row 271, appended at the end of the data block (what used to be row 270)
$row271 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44210, 9, 100112005, 'Puerro', 'Azul de Maquehue', 'Primera', 65, 12000, 12000, 12000, '$/docena de paquetes', 'Provincia de Cautín', 1000, 12, 'Hortaliza')
for ($c = 1; $c -le $row271.Length; $c++) {
    $ws.Cells.Item(271, $c).Value = $row271[$c - 1]
}

# Column D holds a date stored as a serial number; match the date number format
# used by the rest of the column (e.g. row 179) since the new row has no format yet.
$ws.Cells.Item(271, 4).NumberFormat = $ws.Cells.Item(179, 4).NumberFormat
